$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 values
$ws.Range("B6").Value = "RJ TV 1"
$ws.Range("F6").Value = "Mutirão CadÚnico em Campos. Serão distribuídas 800 senhas. Recadastramento acontece no Turf Club. Repórter *ao vivo*. Equipe esteve antes de abrir o portão no local. Mutirão segue até sexta-feira na Fundação de Esportes. Entrevista com beneficiários, reclamando da espera, das filas e da falta de banheiro no local. Um dos entrevistados chegou a fazer apelo ao prefeito. Entrevista com a coordenadora do CadÚnico, Kamila Oliveira. Matéria informativa, mas com muitas reclamações. *sem nota*"

# Remove row 7 entirely (data previously there is now gone)
$ws.Rows("7:7").Delete()
